$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "Aspect of the code" (one run) -> "Aspect of the C" + "ode"
#         (two runs, same rPr; net visible text becomes "Aspect of the Code")
# ---------------------------------------------------------------------
$table = $d.Tables.Item(1)
$headerCell = $table.Cell(1, 1)
$cellRange = $headerCell.Range

# The cell range includes the trailing cell-mark, so grab just the 19
# characters of "Aspect of the code".
$runRange = $d.Range($cellRange.Start, $cellRange.Start + 19)

$splitXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p w:rsidR="000A6CB7" w:rsidRDefault="000A6CB7" w:rsidP="002D2C18">' +
  '<w:pPr><w:spacing w:before="240" w:after="240" w:line="480" w:lineRule="auto"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Aspect of the C</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>ode</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$runRange.InsertXML($splitXml) | Out-Null

# ---------------------------------------------------------------------
# Edit 2: delete the "Table 1: GRNsight ..." caption run entirely,
#         leaving the (now empty) paragraph with its bookmark intact.
# ---------------------------------------------------------------------
$captionText = "Table 1: GRNsight test suite code coverage summary. " +
  "Denominators represent the number of aspects of each type detected " +
  "by Istanbul in the GRNsight codebase; numerators represent the " +
  "subset of these which were executed by unit test code."

$searchRange = $d.Content
$found = $searchRange.Find.Execute($captionText, $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if ($found) {
    $searchRange.Delete()
}
